# Daily auto-push update: insert a new reading row at row 567
# (date 2026/01/06, day 火, time 12, ranking 158), shifting every
# subsequent row down by one (old row 567 -> 568, ..., old row 608 -> 609).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 567; everything from the old
# row 567 onward (through 608) shifts down to 568..609.
$ws.Rows.Item(567).Insert()

# Populate the newly inserted row 567.
# Column A holds a date formatted as plain text ("yyyy/mm/dd"), so force
# the cell to Text format first to stop Excel from auto-converting the
# string into a date serial number; restore the cell style afterward so
# it matches the unstyled data cells around it.
$ws.Range("A567").NumberFormat = "@"
$ws.Range("A567").Value = "2026/01/06"
$ws.Range("A567").Style = "Normal"

$ws.Range("B567").Value = "火"
$ws.Range("C567").Value = 12
$ws.Range("D567").Value = 158
